$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new portfolio data row (row 63) for 2025-10-17.
$dateCell = $ws.Range("A63")

# A leading apostrophe forces Excel to store the date-looking string as
# literal text instead of auto-converting it to a date serial number.
# Resetting the cell style to "Normal" afterwards drops the quote-prefix
# styling so the cell keeps the plain (unstyled) look used by the rest of
# the data rows, while the stored text itself stays the clean
# "2025-10-17" (no apostrophe baked into the value).
$dateCell.Value = "'2025-10-17"
$dateCell.Style = "Normal"

$ws.Range("B63").Value = 52.91999816894531
$ws.Range("C63").Value = 396.6000061035156
$ws.Range("D63").Value = 342.6499938964844
